$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting (it holds
# values like thousands-separated prices and small decimals that Excel would
# otherwise auto-convert to numbers), by pre-formatting it as Text.
$ws.Range("D2:D51").NumberFormat = "@"

# Rows 42 and 43 had their Coin name and Link swapped (Quant now ranks above RenderToken)
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"

# Refresh the Price (D) and Volume(1h) (E) columns for every coin row
$ws.Range("D2").Value = "25.835.01"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").Value = "1.749.90"
$ws.Range("E3").Value = "  -4.95%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "238.32"
$ws.Range("E5").Value = "  -8.12%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.5068"
$ws.Range("E7").Value = "  -5.39%  "
$ws.Range("D8").Value = "41.79"
$ws.Range("E8").Value = "  -6.85%  "
$ws.Range("D9").Value = "0.2702"
$ws.Range("E9").Value = "  -6.33%  "
$ws.Range("D10").Value = "0.06164"
$ws.Range("E10").Value = "  -12.06%  "
$ws.Range("D11").Value = "1.756.31"
$ws.Range("E11").Value = "  -4.66%  "
$ws.Range("D12").Value = "0.06948"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").Value = "15.58"
$ws.Range("E13").Value = "  -8.65%  "
$ws.Range("D14").Value = "0.6021"
$ws.Range("E14").Value = "  -14.72%  "
$ws.Range("D15").Value = "4.494"
$ws.Range("E15").Value = "  -8.90%  "
$ws.Range("D16").Value = "76.83"
$ws.Range("E16").Value = "  -13.24%  "
$ws.Range("D17").Value = "0.9999"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "25.833.30"
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("D20").Value = "0.000006832"
$ws.Range("E20").Value = "  -12.87%  "
$ws.Range("D21").Value = "11.67"
$ws.Range("E21").Value = "  -13.62%  "
$ws.Range("D22").Value = "1.973.44"
$ws.Range("E22").Value = "  -4.95%  "
$ws.Range("D23").Value = "4.052"
$ws.Range("E23").Value = "  -10.76%  "
$ws.Range("D24").Value = "8.196"
$ws.Range("E24").Value = "  -10.12%  "
$ws.Range("D25").Value = "5.208"
$ws.Range("E25").Value = "  -12.19%  "
$ws.Range("D26").Value = "137.15"
$ws.Range("E26").Value = "  -3.89%  "
$ws.Range("D27").Value = "1.466"
$ws.Range("E27").Value = "  -13.81%  "
$ws.Range("D28").Value = "14.98"
$ws.Range("E28").Value = "  -11.13%  "
$ws.Range("D29").Value = "1.803"
$ws.Range("E29").Value = "  -12.86%  "
$ws.Range("D30").Value = "102.59"
$ws.Range("E30").Value = "  -7.35%  "
$ws.Range("D31").Value = "0.08199"
$ws.Range("E31").Value = "  -6.73%  "
$ws.Range("D32").Value = "3.691"
$ws.Range("E32").Value = "  -12.56%  "
$ws.Range("D33").Value = "3.496"
$ws.Range("E33").Value = "  -12.32%  "
$ws.Range("D34").Value = "0.04499"
$ws.Range("E34").Value = "  -5.09%  "
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "2.649"
$ws.Range("E36").Value = "  -8.80%  "
$ws.Range("D37").Value = "0.9833"
$ws.Range("E37").Value = "  -12.59%  "
$ws.Range("D38").Value = "0.6026"
$ws.Range("E38").Value = "  -15.85%  "
$ws.Range("D39").Value = "2.697"
$ws.Range("E39").Value = "  -12.82%  "
$ws.Range("D40").Value = "0.01551"
$ws.Range("E40").Value = "  -7.37%  "
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "103.91"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").Value = "1.889"
$ws.Range("E43").Value = "  -15.11%  "
$ws.Range("D44").Value = "0.3839"
$ws.Range("E44").Value = "  -15.63%  "
$ws.Range("D45").Value = "0.7355"
$ws.Range("E45").Value = "  -17.54%  "
$ws.Range("D46").Value = "4.970"
$ws.Range("E46").Value = "  -14.36%  "
$ws.Range("D47").Value = "0.05425"
$ws.Range("E47").Value = "  -4.40%  "
$ws.Range("D48").Value = "0.1108"
$ws.Range("E48").Value = "  -9.10%  "
$ws.Range("D49").Value = "5.968"
$ws.Range("E49").Value = "  -17.65%  "
$ws.Range("D50").Value = "7.684"
$ws.Range("E50").Value = "  -13.91%  "
$ws.Range("D51").Value = "29.89"
$ws.Range("E51").Value = "  -12.48%  "
